$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Carry the existing date-column formatting (style index with numFmtId 14,
# "m/d/yyyy") down into the two new rows before writing date values, so
# Excel doesn't mint a brand-new custom number format for them.
$ws.Range("C106:D106").Copy()
$ws.Range("C107:D108").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 107: The 80/20 Principle and 92 Other Powerful Laws of Nature
$ws.Cells.Item(107, 1).Value = "The 80/20 Principle and 92 Other Powerful Laws of Nature"
$ws.Cells.Item(107, 2).Value = "Richard Koch"
$ws.Cells.Item(107, 3).Value = (Get-Date -Year 2020 -Month 7 -Day 10 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(107, 4).Value = (Get-Date -Year 2020 -Month 7 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(107, 5).Value = "business;science;success;80/20 rule;darwin;evolution"
$ws.Cells.Item(107, 6).Value = "Audio"
$ws.Cells.Item(107, 7).Value = "12 Hours 10 Mins"
$ws.Cells.Item(107, 8).Value = 3
$ws.Cells.Item(107, 9).Value = $true

# Row 108: The Third Door
$ws.Cells.Item(108, 1).Value = "The Third Door"
$ws.Cells.Item(108, 2).Value = "Alex Banayan"
$ws.Cells.Item(108, 3).Value = (Get-Date -Year 2020 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(108, 4).Value = (Get-Date -Year 2020 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(108, 5).Value = "success;interviewing;failure;tim ferriss;bill gates;jessica alba;larry king;exponential growth"
$ws.Cells.Item(108, 6).Value = "Audio"
$ws.Cells.Item(108, 7).Value = "8 Hours 47 Mins"
$ws.Cells.Item(108, 8).Value = 3
$ws.Cells.Item(108, 9).Value = $true

# Update the view state to match the author's saved selection
$excel.ActiveWindow.ScrollRow = 81
$ws.Range("A109").Select()
